# Update the dSF column (F) values on the active sheet to reflect the
# repulled/pushed data and recalculated mean, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$values = @{
    3  = 6
    4  = -1
    5  = 1
    7  = -1
    8  = 0
    9  = -2
    10 = 1
    11 = 11
    12 = -7
    13 = 1
    14 = -2
    15 = 2
    16 = -1
    17 = -1
    18 = -2
    19 = 4
    20 = -4
    21 = 6
    22 = -1
    23 = -5
    24 = 8
    25 = 2
    26 = -5
    27 = 1
    28 = -1
    30 = 1
    32 = -1
    33 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
